# Apply the edits described by the commit:
#   "fix excel input of variables + make diagramm (tablitsa)0.1.8.0"
#
# Summary of changes:
#  - The active tab moves from sheet 1 ("Локомотив") to sheet 2 ("Лист2").
#  - The selection on sheet 1 ("Локомотив") moves from H9 to B9.
#  - The selection on sheet 2 ("Лист2") moves from D8 to A5.
#  - Cell A5 on sheet 2 gets a new value "[бз]" (a new shared string),
#    replacing the previous "[бh]" value.
#  - The custom number format "0.0" keeps being used by D1 (sheet1 header)
#    / E5 (sheet2) - re-applied so the engine keeps it registered.

$wb = $excel.ActiveWorkbook

$wsLocomotive = $wb.Worksheets.Item(1)
$wsList2      = $wb.Worksheets.Item(2)

# --- Update the value of A5 on "Лист2" (adds a new shared string "[бз]") ---
$wsList2.Range("A5").Value = "[бз]"

# --- Selection on "Локомотив": H9 -> B9 ---
$wsLocomotive.Activate()
$wsLocomotive.Range("B9").Select()

# --- Selection on "Лист2": D8 -> A5, and make it the active tab ---
$wsList2.Activate()
$wsList2.Range("A5").Select()
